$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry pairs a cell reference with its new text value. Every value is
# forced to Text format before being written so numeric-looking strings
# (e.g. "1.00", "0.998", "69.901.23", "  -1.90%  ") are preserved verbatim
# instead of being auto-converted to numbers/percentages by Excel.
$updates = @(
    @{ Cell = 'D2'; Value = '69.901.23' },
    @{ Cell = 'E2'; Value = '  -1.90%  ' },
    @{ Cell = 'D3'; Value = '3.535.32' },
    @{ Cell = 'E3'; Value = '  -1.42%  ' },
    @{ Cell = 'E4'; Value = '  -0.03%  ' },
    @{ Cell = 'D5'; Value = '613.92' },
    @{ Cell = 'E5'; Value = '  +5.10%  ' },
    @{ Cell = 'D6'; Value = '187.23' },
    @{ Cell = 'E6'; Value = '  +0.31%  ' },
    @{ Cell = 'D7'; Value = '0.632' },
    @{ Cell = 'E7'; Value = '  +1.58%  ' },
    @{ Cell = 'D8'; Value = '1.00' },
    @{ Cell = 'E8'; Value = '  -0.10%  ' },
    @{ Cell = 'D9'; Value = '0.216' },
    @{ Cell = 'E9'; Value = '  +0.50%  ' },
    @{ Cell = 'D10'; Value = '0.659' },
    @{ Cell = 'E10'; Value = '  +1.06%  ' },
    @{ Cell = 'D11'; Value = '53.55' },
    @{ Cell = 'E11'; Value = '  -1.78%  ' },
    @{ Cell = 'D12'; Value = '0.0000307' },
    @{ Cell = 'E12'; Value = '  -3.38%  ' },
    @{ Cell = 'D13'; Value = '9.69' },
    @{ Cell = 'E13'; Value = '  +1.59%  ' },
    @{ Cell = 'D14'; Value = '4.107.29' },
    @{ Cell = 'E14'; Value = '  +2.44%  ' },
    @{ Cell = 'D15'; Value = '617.84' },
    @{ Cell = 'E15'; Value = '  +9.44%  ' },
    @{ Cell = 'D16'; Value = '12.84' },
    @{ Cell = 'E16'; Value = '  +4.03%  ' },
    @{ Cell = 'B17'; Value = 'WrappedBTC' },
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc' },
    @{ Cell = 'D17'; Value = '70.005.83' },
    @{ Cell = 'E17'; Value = '  -1.70%  ' },
    @{ Cell = 'B18'; Value = 'Chainlink' },
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link' },
    @{ Cell = 'D18'; Value = '19.14' },
    @{ Cell = 'E18'; Value = '  -0.58%  ' },
    @{ Cell = 'D19'; Value = '3.557.94' },
    @{ Cell = 'E19'; Value = '  -0.41%  ' },
    @{ Cell = 'E20'; Value = '  +0.11%  ' },
    @{ Cell = 'D21'; Value = '0.998' },
    @{ Cell = 'E21'; Value = '  -1.42%  ' },
    @{ Cell = 'D22'; Value = '17.57' },
    @{ Cell = 'E22'; Value = '  -0.06%  ' },
    @{ Cell = 'D23'; Value = '105.44' },
    @{ Cell = 'E23'; Value = '  +11.15%  ' },
    @{ Cell = 'D24'; Value = '4.69' },
    @{ Cell = 'E24'; Value = '  +2.42%  ' },
    @{ Cell = 'D25'; Value = '5.04' },
    @{ Cell = 'E25'; Value = '  -1.35%  ' },
    @{ Cell = 'D26'; Value = '3.03' },
    @{ Cell = 'E26'; Value = '  +3.07%  ' },
    @{ Cell = 'D27'; Value = '10.93' },
    @{ Cell = 'E27'; Value = '  -2.96%  ' },
    @{ Cell = 'D28'; Value = '10.06' },
    @{ Cell = 'E28'; Value = '  +9.83%  ' },
    @{ Cell = 'E29'; Value = '  +5.07%  ' },
    @{ Cell = 'D30'; Value = '7.05' },
    @{ Cell = 'E30'; Value = '  -2.89%  ' },
    @{ Cell = 'D31'; Value = '12.47' },
    @{ Cell = 'E31'; Value = '  +1.37%  ' },
    @{ Cell = 'D32'; Value = '0.116' },
    @{ Cell = 'E32'; Value = '  +0.98%  ' },
    @{ Cell = 'D33'; Value = '63.95' },
    @{ Cell = 'E33'; Value = '  -0.52%  ' },
    @{ Cell = 'D34'; Value = '3.73' },
    @{ Cell = 'E34'; Value = '  +12.39%  ' },
    @{ Cell = 'D35'; Value = '536.71' },
    @{ Cell = 'E35'; Value = '  -1.79%  ' },
    @{ Cell = 'B36'; Value = 'Dai' },
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai' },
    @{ Cell = 'D36'; Value = '1.00' },
    @{ Cell = 'E36'; Value = '  -0.02%  ' },
    @{ Cell = 'B37'; Value = 'Fetch.AI' },
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet' },
    @{ Cell = 'D37'; Value = '3.13' },
    @{ Cell = 'E37'; Value = '  -6.48%  ' },
    @{ Cell = 'D38'; Value = '0.399' },
    @{ Cell = 'E38'; Value = '  -4.36%  ' },
    @{ Cell = 'B39'; Value = 'Stacks' },
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx' },
    @{ Cell = 'D39'; Value = '3.59' },
    @{ Cell = 'E39'; Value = '  +3.55%  ' },
    @{ Cell = 'B40'; Value = 'InjectiveProtocol' },
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj' },
    @{ Cell = 'D40'; Value = '36.82' },
    @{ Cell = 'E40'; Value = '  -2.16%  ' },
    @{ Cell = 'B41'; Value = 'Maker' },
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr' },
    @{ Cell = 'D41'; Value = '3.548.78' },
    @{ Cell = 'E41'; Value = '  +1.14%  ' },
    @{ Cell = 'B42'; Value = 'PEPE' },
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe' },
    @{ Cell = 'D42'; Value = '0.0₃0777' },
    @{ Cell = 'E42'; Value = '  -3.71%  ' },
    @{ Cell = 'E43'; Value = '  +3.31%  ' },
    @{ Cell = 'D44'; Value = '0.0463' },
    @{ Cell = 'E44'; Value = '  +3.93%  ' },
    @{ Cell = 'D45'; Value = '2.96' },
    @{ Cell = 'E45'; Value = '  +0.74%  ' },
    @{ Cell = 'E46'; Value = '  +4.67%  ' },
    @{ Cell = 'D47'; Value = '3.36' },
    @{ Cell = 'E47'; Value = '  -3.37%  ' },
    @{ Cell = 'D48'; Value = '8.96' },
    @{ Cell = 'E48'; Value = '  -4.38%  ' },
    @{ Cell = 'E49'; Value = '  +0.30%  ' },
    @{ Cell = 'D50'; Value = '132.31' },
    @{ Cell = 'E50'; Value = '  -1.90%  ' },
    @{ Cell = 'D51'; Value = '1.36' },
    @{ Cell = 'E51'; Value = '  -6.03%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}

